$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.105.04"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "1.640.06"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "1.868.47"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "1.623.38"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "26.112.95"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("E35").Value = "  +1.37%  "

$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("D37").Value = "1.134.85"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.551"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.796"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "1.777.89"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.33%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.416"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
